$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, shifting N:P -> O:Q
$ws.Columns("N:N").Insert()

# New column N inherits the column width of the column to its left (M),
# matching Excel's native "insert column" behavior
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab
$ws.Activate()

# Update the selected cell on the Repayment schedule sheet
$ws.Range("S4").Select()
